$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("price1")
$ws.Range("A2").Value = "price1 Step5 Seed"
$ws.Range("A3").Value = "price1 Step4 Downgrade"
$ws.Range("A4").Value = "price1 Step3 Challenge"
$ws.Range("A5").Value = "price1 Step2 FOMO"
$ws.Range("A6").Value = "price1 Step1 Reframe"

$ws = $wb.Worksheets.Item("price2")
$ws.Range("A2").Value = "price2 Step5 Seed"
$ws.Range("A3").Value = "price2 Step4 Downgrade"
$ws.Range("A4").Value = "price2 Step3 Challenge"
$ws.Range("A5").Value = "price2 Step2 FOMO"
$ws.Range("A6").Value = "price2 Step1 Reframe"

$ws = $wb.Worksheets.Item("discount1")
$ws.Range("A2").Value = "discount1 Step4 Takeaway"
$ws.Range("A3").Value = "discount1 Step3 Concession"
$ws.Range("A4").Value = "discount1 Step2 Challenge"
$ws.Range("A5").Value = "discount1 Step1 Firmness"

$ws = $wb.Worksheets.Item("discount2")
$ws.Range("A2").Value = "discount2 Step4 Takeaway"
$ws.Range("A3").Value = "discount2 Step3 Concession"
$ws.Range("A4").Value = "discount2 Step2 Challenge"
$ws.Range("A5").Value = "discount2 Step1 Firmness"

$ws = $wb.Worksheets.Item("free1")
$ws.Range("A2").Value = "free1 Step4 Seed"
$ws.Range("A3").Value = "free1 Step3 Guilt"
$ws.Range("A4").Value = "free1 Step2 Challenge"
$ws.Range("A5").Value = "free1 Step1 Reminder"

$ws = $wb.Worksheets.Item("free2")
$ws.Range("A2").Value = "free2 Step4 Seed"
$ws.Range("A3").Value = "free2 Step3 Guilt"
$ws.Range("A4").Value = "free2 Step2 Challenge"
$ws.Range("A5").Value = "free2 Step1 Reminder"

$ws = $wb.Worksheets.Item("nomoney1")
$ws.Range("A2").Value = "nomoney1 Step4 Protect"
$ws.Range("A3").Value = "nomoney1 Step3 PWYW"
$ws.Range("A4").Value = "nomoney1 Step2 Test"
$ws.Range("A5").Value = "nomoney1 Step1 Empathy"

$ws = $wb.Worksheets.Item("nomoney2")
$ws.Range("A2").Value = "nomoney2 Step4 Protect"
$ws.Range("A3").Value = "nomoney2 Step3 PWYW"
$ws.Range("A4").Value = "nomoney2 Step2 Test"
$ws.Range("A5").Value = "nomoney2 Step1 Empathy"

$ws = $wb.Worksheets.Item("noppv1")
$ws.Range("A2").Value = "noppv1 Step3 PWYW"
$ws.Range("A3").Value = "noppv1 Step2 Reframe"
$ws.Range("A4").Value = "noppv1 Step1 Accept"

$ws = $wb.Worksheets.Item("noppv2")
$ws.Range("A2").Value = "noppv2 Step3 PWYW"
$ws.Range("A3").Value = "noppv2 Step2 Reframe"
$ws.Range("A4").Value = "noppv2 Step1 Accept"

$ws = $wb.Worksheets.Item("card1")
$ws.Range("A2").Value = "card1 Step3 Urgency"
$ws.Range("A3").Value = "card1 Step2 AltCard"
$ws.Range("A4").Value = "card1 Step1 Retry"

$ws = $wb.Worksheets.Item("card2")
$ws.Range("A2").Value = "card2 Step3 Urgency"
$ws.Range("A3").Value = "card2 Step2 AltCard"
$ws.Range("A4").Value = "card2 Step1 Retry"

$ws = $wb.Worksheets.Item("nosex1")
$ws.Range("A2").Value = "nosex1 Step4 Accept"
$ws.Range("A3").Value = "nosex1 Step3 ReAttempt"
$ws.Range("A4").Value = "nosex1 Step2 Subtle"
$ws.Range("A5").Value = "nosex1 Step1 Respect"

$ws = $wb.Worksheets.Item("nosex2")
$ws.Range("A2").Value = "nosex2 Step4 Accept"
$ws.Range("A3").Value = "nosex2 Step3 ReAttempt"
$ws.Range("A4").Value = "nosex2 Step2 Subtle"
$ws.Range("A5").Value = "nosex2 Step1 Respect"

$ws = $wb.Worksheets.Item("offtopic1")
$ws.Range("A2").Value = "offtopic1 Step3 Retake"
$ws.Range("A3").Value = "offtopic1 Step2 Redirect"
$ws.Range("A4").Value = "offtopic1 Step1 Acknowledge"

$ws = $wb.Worksheets.Item("offtopic2")
$ws.Range("A2").Value = "offtopic2 Step3 Retake"
$ws.Range("A3").Value = "offtopic2 Step2 Redirect"
$ws.Range("A4").Value = "offtopic2 Step1 Acknowledge"

$ws = $wb.Worksheets.Item("real1")
$ws.Range("A2").Value = "real1 Step3 Grounding"
$ws.Range("A3").Value = "real1 Step2 Challenge"
$ws.Range("A4").Value = "real1 Step1 Humor"

$ws = $wb.Worksheets.Item("real2")
$ws.Range("A2").Value = "real2 Step3 Grounding"
$ws.Range("A3").Value = "real2 Step2 Challenge"
$ws.Range("A4").Value = "real2 Step1 Humor"

$ws = $wb.Worksheets.Item("voice1")
$ws.Range("A2").Value = "voice1 Step3 Firm"
$ws.Range("A3").Value = "voice1 Step2 Redirect"
$ws.Range("A4").Value = "voice1 Step1 Dodge"

$ws = $wb.Worksheets.Item("voice2")
$ws.Range("A2").Value = "voice2 Step3 Firm"
$ws.Range("A3").Value = "voice2 Step2 Redirect"
$ws.Range("A4").Value = "voice2 Step1 Dodge"

$ws = $wb.Worksheets.Item("customyes1")
$ws.Range("A2").Value = "customyes1 Step3 Close"
$ws.Range("A3").Value = "customyes1 Step2 Price"
$ws.Range("A4").Value = "customyes1 Step1 Tease"

$ws = $wb.Worksheets.Item("customyes2")
$ws.Range("A2").Value = "customyes2 Step3 Close"
$ws.Range("A3").Value = "customyes2 Step2 Price"
$ws.Range("A4").Value = "customyes2 Step1 Tease"

$ws = $wb.Worksheets.Item("customno1")
$ws.Range("A2").Value = "customno1 Step3 Close"
$ws.Range("A3").Value = "customno1 Step2 Alternative"
$ws.Range("A4").Value = "customno1 Step1 Redirect"

$ws = $wb.Worksheets.Item("customno2")
$ws.Range("A2").Value = "customno2 Step3 Close"
$ws.Range("A3").Value = "customno2 Step2 Alternative"
$ws.Range("A4").Value = "customno2 Step1 Redirect"

$ws = $wb.Worksheets.Item("done1")
$ws.Range("A2").Value = "done1 Step3 Seed"
$ws.Range("A3").Value = "done1 Step2 Rescue"
$ws.Range("A4").Value = "done1 Step1 Validate"

$ws = $wb.Worksheets.Item("done2")
$ws.Range("A2").Value = "done2 Step3 Seed"
$ws.Range("A3").Value = "done2 Step2 Rescue"
$ws.Range("A4").Value = "done2 Step1 Validate"

$ws = $wb.Worksheets.Item("cumcontrol")
$ws.Range("A2").Value = "cumcontrol delay2"
$ws.Range("A3").Value = "cumcontrol delay1"
$ws.Range("A4").Value = "cumcontrol sync2"
$ws.Range("A5").Value = "cumcontrol sync1"
$ws.Range("A6").Value = "cumcontrol edge2"
$ws.Range("A7").Value = "cumcontrol edge1"

$ws = $wb.Worksheets.Item("dickpic")
$ws.Range("A2").Value = "dickpic dpppv2"
$ws.Range("A3").Value = "dickpic dpppv1"
$ws.Range("A4").Value = "dickpic dprapport2"
$ws.Range("A5").Value = "dickpic dprapport1"
$ws.Range("A6").Value = "dickpic dpsext2"
$ws.Range("A7").Value = "dickpic dpsext1"

$ws = $wb.Worksheets.Item("boosters")
$ws.Range("A2").Value = "boosters h8"
$ws.Range("A3").Value = "boosters h7"
$ws.Range("A4").Value = "boosters h6"
$ws.Range("A5").Value = "boosters h5"
$ws.Range("A6").Value = "boosters h4"
$ws.Range("A7").Value = "boosters h3"
$ws.Range("A8").Value = "boosters h2"
$ws.Range("A9").Value = "boosters h1"
